{"js": "// Apply the git.docx edit:\n//  1) Append a new run \" or git add . \\u2013 this will stage all the changes\"\n//     to the end of the \"Git add <filename>- to stage the file\" paragraph.\n//  2) Add a new paragraph \"Git init- initializes repository\" right after\n//     the \"Git log ...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet addParagraph = null;\nlet logParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (addParagraph === null && text.indexOf(\"Git add <filename>- to stage the file\") !== -1) {\n    addParagraph = paragraphs.items[i];\n  }\n  if (logParagraph === null && text.indexOf(\"Git log\") !== -1 && text.indexOf(\"history if all the commits\") !== -1) {\n    logParagraph = paragraphs.items[i];\n  }\n}\n\n// 1) Append the extra sentence as its own run (matches the OOXML diff,\n// which shows a brand-new <w:r> rather than text merged into the\n// existing run). insertOoxml lets us control the run boundary exactly;\n// insertText() would just extend the existing run's text.\nif (addParagraph) {\n  const endRange = addParagraph.getRange(Word.RangeLocation.end);\n  const flatOpcXml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:t xml:space=\"preserve\"> or git add . \\u2013 this will stage all the changes</w:t></w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n  endRange.insertOoxml(flatOpcXml, Word.InsertLocation.end);\n  await context.sync();\n}\n\n// 2) Insert the new \"Git init\" paragraph right after the \"Git log\" paragraph.\nif (logParagraph) {\n  logParagraph.insertParagraph(\"Git init- initializes repository\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Apply the git.docx edit:\n#  1) Append a new run \" or git add . - this will stage all the changes\"\n#     to the end of the \"Git add <filename>- to stage the file\" paragraph.\n#  2) Add a new paragraph \"Git init- initializes repository\" right after\n#     the \"Git log ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n$addParagraph = $null\n$logParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($addParagraph -eq $null -and $t -like \"*Git add <filename>- to stage the file*\") {\n        $addParagraph = $p\n    }\n    if ($logParagraph -eq $null -and $t -like \"*Git log*\" -and $t -like \"*history if all the commits*\") {\n        $logParagraph = $p\n    }\n}\n\n# 1) Append the extra sentence as its own run (matches the OOXML diff,\n# which shows a brand-new <w:r> rather than text merged into the\n# existing run). Range.InsertXML with an explicit \"End\" location keeps\n# the insertion inside the paragraph (instead of splitting it into a new\n# paragraph) while still creating a fresh run boundary.\nif ($addParagraph -ne $null) {\n    $r = $addParagraph.Range\n    $r.MoveEnd(1, -1)  # exclude the paragraph mark\n    $r.Collapse(0)     # wdCollapseEnd -> collapse to a point right after \"file\"\n    $flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:t xml:space=\"preserve\"> or git add . ' + [char]0x2013 + ' this will stage all the changes</w:t></w:r></w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n    $r.InsertXML($flatOpc, \"End\")\n}\n\n# 2) Insert the new \"Git init\" paragraph right after the \"Git log\" paragraph.\nif ($logParagraph -ne $null) {\n    $r2 = $logParagraph.Range\n    $r2.Collapse(0)  # wdCollapseEnd -> point right after the paragraph mark\n    $r2.InsertParagraphAfter()\n    $newParagraph = $logParagraph.Next()\n    $newParagraph.Range.Text = \"Git init- initializes repository\"\n}\n"}
